# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows (for a new report date, 2023-11-28) into the
# "Uva" price list, just above the block of rows that begin at the old row
# 164. Excel's row insert shifts all of the old rows 164-193 down to 167-196,
# preserving their existing values/formatting, and we then populate the 3
# freshly inserted rows with the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows right before the current row 164 (formatting copies
# down from the row above, same as interactive Excel row insert).
$ws.Range("A164:A166").EntireRow.Insert()

# --- Row 164: Red Globe / Tercera, Provincia de Copiapó --------------------
$ws.Cells.Item(164, 1).Value = 1
$ws.Cells.Item(164, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(164, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(164, 4).Value = 45258
$ws.Cells.Item(164, 5).Value = 15
$ws.Cells.Item(164, 6).Value = "Fruta"
$ws.Cells.Item(164, 7).Value = 100109
$ws.Cells.Item(164, 8).Value = "Uva"
$ws.Cells.Item(164, 9).Value = 100109001
$ws.Cells.Item(164, 10).Value = "Uva"
$ws.Cells.Item(164, 11).Value = "Red Globe"
$ws.Cells.Item(164, 12).Value = "Tercera"
$ws.Cells.Item(164, 13).Value = 270
$ws.Cells.Item(164, 14).Value = 24000
$ws.Cells.Item(164, 15).Value = 25000
$ws.Cells.Item(164, 16).Value = 24500
$ws.Cells.Item(164, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(164, 18).Value = "Provincia de Copiapó"
$ws.Cells.Item(164, 19).Value = 2450
$ws.Cells.Item(164, 20).Value = 10

# --- Row 165: Superior Seedless / Primera, Región de Coquimbo --------------
$ws.Cells.Item(165, 1).Value = 1
$ws.Cells.Item(165, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(165, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(165, 4).Value = 45258
$ws.Cells.Item(165, 5).Value = 15
$ws.Cells.Item(165, 6).Value = "Fruta"
$ws.Cells.Item(165, 7).Value = 100109
$ws.Cells.Item(165, 8).Value = "Uva"
$ws.Cells.Item(165, 9).Value = 100109001
$ws.Cells.Item(165, 10).Value = "Uva"
$ws.Cells.Item(165, 11).Value = "Superior Seedless"
$ws.Cells.Item(165, 12).Value = "Primera"
$ws.Cells.Item(165, 13).Value = 200
$ws.Cells.Item(165, 14).Value = 24000
$ws.Cells.Item(165, 15).Value = 25000
$ws.Cells.Item(165, 16).Value = 24500
$ws.Cells.Item(165, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(165, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(165, 19).Value = 2450
$ws.Cells.Item(165, 20).Value = 10

# --- Row 166: Superior Seedless / Segunda, Provincia de Copiapó ------------
$ws.Cells.Item(166, 1).Value = 1
$ws.Cells.Item(166, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(166, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(166, 4).Value = 45258
$ws.Cells.Item(166, 5).Value = 15
$ws.Cells.Item(166, 6).Value = "Fruta"
$ws.Cells.Item(166, 7).Value = 100109
$ws.Cells.Item(166, 8).Value = "Uva"
$ws.Cells.Item(166, 9).Value = 100109001
$ws.Cells.Item(166, 10).Value = "Uva"
$ws.Cells.Item(166, 11).Value = "Superior Seedless"
$ws.Cells.Item(166, 12).Value = "Segunda"
$ws.Cells.Item(166, 13).Value = 270
$ws.Cells.Item(166, 14).Value = 24000
$ws.Cells.Item(166, 15).Value = 25000
$ws.Cells.Item(166, 16).Value = 24500
$ws.Cells.Item(166, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(166, 18).Value = "Provincia de Copiapó"
$ws.Cells.Item(166, 19).Value = 2450
$ws.Cells.Item(166, 20).Value = 10
